# Update countries & provincias Spain
# This applies updated COVID-19 case numbers for several countries and
# re-sorts the "Maldivas"/"Congo" rows (table is sorted descending by
# "Casos totales", column B) since Maldivas overtook Congo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos (Estados Unidos keeps first place, only totals updated)
$ws.Range("B4").Value = 980784
$ws.Range("C4").Value = 20133
$ws.Range("D4").Value = 118735
$ws.Range("E4").Value = 806955
$ws.Range("G4").Value = 838
$ws.Range("H4").Value = 55094

# Row 8: Alemania
$ws.Range("B8").Value = 157177
$ws.Range("C8").Value = 664
$ws.Range("E8").Value = 39264
$ws.Range("G8").Value = 36
$ws.Range("H8").Value = 5913

# Rows 130-131: Maldivas overtakes Congo (214 > 200 casos totales), so the
# row order is swapped to keep the sheet sorted by Casos totales (col B).
$ws.Range("A130").Value = "Maldivas"
$ws.Range("B130").Value = 214
$ws.Range("C130").Value = 37
$ws.Range("D130").Value = 17
$ws.Range("E130").Value = 197
$ws.Range("F130").Value = 2
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0

$ws.Range("A131").Value = "Congo"
$ws.Range("B131").Value = 200
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 19
$ws.Range("E131").Value = 175
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 6

# Row 149: Togo
$ws.Range("B149").Value = 98
$ws.Range("C149").Value = 2
$ws.Range("E149").Value = 30
